# Rotate the data of rows 10, 11 and 14 on Sheet1.
# Before:
#   Row 10: Jarrett Allen        | C    | Cleveland Cavaliers
#   Row 11: Karl-Anthony Towns   | PF,C | New York Knicks
#   Row 14: Daniel Gafford       | PF,C | Dallas Mavericks
# After:
#   Row 10: Karl-Anthony Towns   | PF,C | New York Knicks
#   Row 11: Daniel Gafford       | PF,C | Dallas Mavericks
#   Row 14: Jarrett Allen        | C    | Cleveland Cavaliers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 gets what used to be in row 11 (Karl-Anthony Towns).
$ws.Range("A10").Value = "Karl-Anthony Towns"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "New York Knicks"

# Row 11 gets what used to be in row 14 (Daniel Gafford).
$ws.Range("A11").Value = "Daniel Gafford"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Dallas Mavericks"

# Row 14 gets what used to be in row 10 (Jarrett Allen).
$ws.Range("A14").Value = "Jarrett Allen"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Cleveland Cavaliers"
